$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "ECs sender" rows (old rows 8-10); remaining rows shift up
$ws.Rows("8:10").Delete()

# Rewrite the updated TPM-derived values for rows 2-7
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf16"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7297526666666666
$ws.Range("H2").Value = 2.189258
$ws.Range("I2").Value = 0.2331354772809149
$ws.Range("J2").Value = 0.2331354772809149
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.045145666666667
$ws.Range("N2").Value = 18.135437
$ws.Range("O2").Value = 0.8160840232643366
$ws.Range("P2").Value = 0.8160840232643367
$ws.Range("Q2").Value = 4.411461170638445
$ws.Range("R2").Value = 39.703150535746
$ws.Range("S2").Value = 0.1902581382650604
$ws.Range("T2").Value = 0.1902581382650604

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf16"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.7297526666666666
$ws.Range("H3").Value = 2.189258
$ws.Range("I3").Value = 0.2331354772809149
$ws.Range("J3").Value = 0.2331354772809149
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6824433333333334
$ws.Range("N3").Value = 2.04733
$ws.Range("O3").Value = 0.09212864864242169
$ws.Range("P3").Value = 0.09212864864242169
$ws.Range("Q3").Value = 0.4980148423488889
$ws.Range("R3").Value = 4.482133581139999
$ws.Range("S3").Value = 0.0214784564724967
$ws.Range("T3").Value = 0.02147845647249669

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf16"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.7297526666666666
$ws.Range("H4").Value = 2.189258
$ws.Range("I4").Value = 0.2331354772809149
$ws.Range("J4").Value = 0.2331354772809149
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6799149999999999
$ws.Range("N4").Value = 2.039745
$ws.Range("O4").Value = 0.09178732809324164
$ws.Range("P4").Value = 0.09178732809324165
$ws.Range("Q4").Value = 0.4961697843566666
$ws.Range("R4").Value = 4.46552805921
$ws.Range("S4").Value = 0.02139888254335782
$ws.Range("T4").Value = 0.02139888254335782

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Fgf16"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.400413
$ws.Range("H5").Value = 7.201238999999999
$ws.Range("I5").Value = 0.766864522719085
$ws.Range("J5").Value = 0.766864522719085
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.045145666666667
$ws.Range("N5").Value = 18.135437
$ws.Range("O5").Value = 0.8160840232643366
$ws.Range("P5").Value = 0.8160840232643367
$ws.Range("Q5").Value = 14.51084624516033
$ws.Range("R5").Value = 130.597616206443
$ws.Range("S5").Value = 0.6258258849992762
$ws.Range("T5").Value = 0.6258258849992763

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf16"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.400413
$ws.Range("H6").Value = 7.201238999999999
$ws.Range("I6").Value = 0.766864522719085
$ws.Range("J6").Value = 0.766864522719085
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.6824433333333334
$ws.Range("N6").Value = 2.04733
$ws.Range("O6").Value = 0.09212864864242169
$ws.Range("P6").Value = 0.09212864864242169
$ws.Range("Q6").Value = 1.638145849096667
$ws.Range("R6").Value = 14.74331264187
$ws.Range("S6").Value = 0.07065019216992499
$ws.Range("T6").Value = 0.07065019216992499

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf16"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.400413
$ws.Range("H7").Value = 7.201238999999999
$ws.Range("I7").Value = 0.766864522719085
$ws.Range("J7").Value = 0.766864522719085
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6799149999999999
$ws.Range("N7").Value = 2.039745
$ws.Range("O7").Value = 0.09178732809324164
$ws.Range("P7").Value = 0.09178732809324165
$ws.Range("Q7").Value = 1.632076804895
$ws.Range("R7").Value = 14.688691244055
$ws.Range("S7").Value = 0.07038844554988381
$ws.Range("T7").Value = 0.07038844554988383
